$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number & report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# --- Cells changing from a number to the "N/A" text marker (copy style+text from donor, style 13) ---
# Donor C14 already holds text "0" with style 13
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("C14").Copy($ws.Range("D28"))

# --- Cells changing from a number to the "***.*" text marker (copy style+text from donor, style 13) ---
# Donor E14 already holds text "***.*" with style 13
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("E14").Copy($ws.Range("H27"))
$ws.Range("E14").Copy($ws.Range("E28"))

# --- Cells changing from "N/A" text to a number (copy number style 14 from donor J14, then set value) ---
$ws.Range("J14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3
$ws.Range("J14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 4
$ws.Range("J14").Copy($ws.Range("C31"))
$ws.Range("C31").Value = 1
$ws.Range("J14").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 2
$ws.Range("J14").Copy($ws.Range("I31"))
$ws.Range("I31").Value = 2

# --- Remaining same-type value updates ---
# Row 15
$ws.Range("M15").Value = 0
# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 18.181818181818
$ws.Range("L16").Value = 36.842105263157
$ws.Range("M16").Value = -46.938775510204
$ws.Range("N16").Value = -86.802030456852
# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -23.076923076923
$ws.Range("I17").Value = 35
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -25.531914893617
$ws.Range("L17").Value = -30
$ws.Range("M17").Value = 9.375
$ws.Range("N17").Value = -44.444444444444
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 31
$ws.Range("K18").Value = 16.129032258064
$ws.Range("L18").Value = 9.090909090909
$ws.Range("M18").Value = -60
$ws.Range("N18").Value = -91.762013729977
# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 166.666666666667
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 13.888888888888
$ws.Range("I19").Value = 105
$ws.Range("J19").Value = 139
$ws.Range("K19").Value = -24.460431654676
$ws.Range("L19").Value = -23.91304347826
$ws.Range("M19").Value = 15.384615384615
$ws.Range("N19").Value = -33.12101910828
# Row 20
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 450
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 127.272727272727
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = 47.058823529411
$ws.Range("M20").Value = 38.888888888888
$ws.Range("N20").Value = -89.339019189765
# Row 21
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 72.727272727272
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 16.883116883116
$ws.Range("I21").Value = 254
$ws.Range("J21").Value = 288
$ws.Range("K21").Value = -11.805555555555
$ws.Range("L21").Value = -8.633093525179
$ws.Range("M21").Value = -15.333333333333
$ws.Range("N21").Value = -80.902255639097
# Row 22
$ws.Range("L22").Value = -20
# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -44
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 5.479452054794
$ws.Range("I24").Value = 269
$ws.Range("J24").Value = 241
$ws.Range("K24").Value = 11.61825726141
$ws.Range("L24").Value = -2.536231884057
$ws.Range("M24").Value = 31.219512195122
# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -5.555555555555
$ws.Range("I25").Value = 49
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -18.333333333333
$ws.Range("L25").Value = -42.35294117647
# Row 26
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 3.225806451612
$ws.Range("I26").Value = 94
$ws.Range("J26").Value = 104
$ws.Range("K26").Value = -9.615384615384
$ws.Range("L26").Value = 34.285714285714
$ws.Range("M26").Value = -6.930693069306
# Row 28
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("K28").Value = -28.571428571428
$ws.Range("L28").Value = 87.5
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("K31").Value = -33.333333333333
$ws.Range("L31").Value = 100
# Row 33
$ws.Range("G33").Value = 2
$ws.Range("J33").Value = 2
